# Updated cryptos list with GitHub Actions: refresh Price/Volume(1h) columns
# and fix the ARBITRUM / WEMIXToken row ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns (rows 2-3 and 5-51) so numeric-looking
# strings like "302.96" or "0.800" are preserved as text, not converted to numbers.
$ws.Range("D2:E3").NumberFormat = "@"
$ws.Range("D5:E51").NumberFormat = "@"

$ws.Range("D2").Value = '43.133.94'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = '2.339.94'
$ws.Range("E3").Value = '  +1.21%  '

$ws.Range("D5").Value = '302.96'
$ws.Range("E5").Value = '  +0.40%  '

$ws.Range("D6").Value = '94.98'
$ws.Range("E6").Value = '  -1.80%  '

$ws.Range("D7").Value = '0.504'
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").Value = '0.496'
$ws.Range("E9").Value = '  -0.30%  '

$ws.Range("D10").Value = '34.13'
$ws.Range("E10").Value = '  -2.47%  '

$ws.Range("D11").Value = '0.0784'
$ws.Range("E11").Value = '  -0.82%  '

$ws.Range("D12").Value = '18.73'
$ws.Range("E12").Value = '  -3.03%  '

$ws.Range("E13").Value = '  +1.89%  '

$ws.Range("D14").Value = '6.75'
$ws.Range("E14").Value = '  -1.28%  '

$ws.Range("D15").Value = '2.701.09'
$ws.Range("E15").Value = '  +0.99%  '

$ws.Range("D16").Value = '2.322.29'
$ws.Range("E16").Value = '  +0.62%  '

$ws.Range("D17").Value = '0.800'
$ws.Range("E17").Value = '  +1.92%  '

$ws.Range("D18").Value = '43.051.77'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").Value = '12.18'
$ws.Range("E19").Value = '  -2.08%  '

$ws.Range("D20").Value = '6.22'
$ws.Range("E20").Value = '  +3.29%  '

$ws.Range("D21").Value = '0.0₃0890'
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").Value = '67.99'
$ws.Range("E22").Value = '  +0.52%  '

$ws.Range("D23").Value = '236.12'
$ws.Range("E23").Value = '  +0.26%  '

$ws.Range("D24").Value = '2.22'
$ws.Range("E24").Value = '  -0.65%  '

$ws.Range("E25").Value = '  +0.22%  '

$ws.Range("E26").Value = '  -0.59%  '

$ws.Range("D27").Value = '24.64'
$ws.Range("E27").Value = '  -0.24%  '

$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  -0.77%  '

$ws.Range("E29").Value = '  +1.20%  '

$ws.Range("D30").Value = '31.56'
$ws.Range("E30").Value = '  -3.06%  '

$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("D32").Value = '5.01'
$ws.Range("E32").Value = '  +0.84%  '

$ws.Range("E33").Value = '  +4.29%  '

$ws.Range("D34").Value = '17.28'
$ws.Range("E34").Value = '  -2.72%  '

$ws.Range("D35").Value = '4.40'
$ws.Range("E35").Value = '  -1.99%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.83'
$ws.Range("E36").Value = '  +3.73%  '

$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D37").Value = '2.32'
$ws.Range("E37").Value = '  -1.04%  '

$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("D39").Value = '2.76'
$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("D40").Value = '22.15'
$ws.Range("E40").Value = '  +18.72%  '

$ws.Range("E41").Value = '  -0.41%  '

$ws.Range("D42").Value = '113.99'
$ws.Range("E42").Value = '  -30.43%  '

$ws.Range("D43").Value = '1.934.74'
$ws.Range("E43").Value = '  -2.20%  '

$ws.Range("D44").Value = '0.0281'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").Value = '10.06'
$ws.Range("E45").Value = '  -5.02%  '

$ws.Range("E46").Value = '  +1.38%  '

$ws.Range("D47").Value = '2.73'
$ws.Range("E47").Value = '  -1.13%  '

$ws.Range("D48").Value = '2.568.56'
$ws.Range("E48").Value = '  +1.17%  '

$ws.Range("E49").Value = '  -0.81%  '

$ws.Range("D50").Value = '53.14'
$ws.Range("E50").Value = '  -0.23%  '

$ws.Range("D51").Value = '72.12'
$ws.Range("E51").Value = '  +0.14%  '
